# Normalize the "Recorded By" (column G) name ordering on the active sheet.
# For each used row, the two comma-separated names in column G are
# reordered so that "System" (when present) is moved to the end of the
# list, and the pairing "dnasr281@gmail.com, admin@admin.com" is swapped
# to "admin@admin.com, dnasr281@gmail.com" - matching the sync performed
# upstream on the attendance report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text.Split(",")
    if ($parts.Count -ne 2) { continue }

    $first = $parts[0].Trim()
    $second = $parts[1].Trim()

    if ($first -eq "System" -or $first -eq "system") {
        $cell.Value = $second + ", " + $first
    } elseif ($first -eq "dnasr281@gmail.com" -and $second -eq "admin@admin.com") {
        $cell.Value = $second + ", " + $first
    }
}
